{"js": "// The document has two \"numbered\" tables whose first (small) column(s)\n// hold a static \"N.\" label and whose last column holds the publication /\n// conference citation text. This edit removes specific entries from each\n// table. Because the \"N.\" labels are plain static text (not an automatic\n// Word numbering field), removing an entry in the middle is implemented by\n// shifting the citation text of every later row up by one slot and then\n// physically deleting the now-duplicated row(s) at the end of the table -\n// this keeps the \"1.\", \"2.\", \"3.\" ... labels correct without having to\n// rewrite them.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Table 1: \"Papers Published by Faculty in Journals\" - 2 columns\n// (label, citation). The first two entries are removed, so the citation\n// text that used to belong to rows 3/4/5 becomes the new rows 1/2/3, and\n// the trailing two rows are deleted.\n// ---------------------------------------------------------------------\nconst journalTable = tables.items[0];\njournalTable.rows.load(\"items\");\nawait context.sync();\n\nconst journalRows = journalTable.rows.items;\nfor (let r = 0; r < journalRows.length; r++) {\n  journalRows[r].cells.load(\"items\");\n}\nawait context.sync();\n\nfor (let r = 0; r < journalRows.length; r++) {\n  const cells = journalRows[r].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    cells[c].body.paragraphs.load(\"items\");\n  }\n}\nawait context.sync();\n\n// Grab the citation text (last cell) of every row before we start\n// overwriting anything.\nconst journalCitations = [];\nfor (let r = 0; r < journalRows.length; r++) {\n  const cells = journalRows[r].cells.items;\n  const lastCellParas = cells[cells.length - 1].body.paragraphs.items;\n  lastCellParas[0].load(\"text\");\n}\nawait context.sync();\nfor (let r = 0; r < journalRows.length; r++) {\n  const cells = journalRows[r].cells.items;\n  const lastCellParas = cells[cells.length - 1].body.paragraphs.items;\n  journalCitations.push(lastCellParas[0].text);\n}\n\nconst journalEntriesRemoved = 2; // \"Tirtharaj ... \" and \"Ashwin ... Shroff ...\" removed\n\nfor (let r = 0; r < journalRows.length - journalEntriesRemoved; r++) {\n  const cells = journalRows[r].cells.items;\n  const lastCellParas = cells[cells.length - 1].body.paragraphs.items;\n  lastCellParas[0].insertText(journalCitations[r + journalEntriesRemoved], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Delete the trailing rows that are now duplicates of the shifted content.\nfor (let i = 0; i < journalEntriesRemoved; i++) {\n  journalTable.rows.load(\"items\");\n  await context.sync();\n  const rows = journalTable.rows.items;\n  rows[rows.length - 1].delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Table 2: \"CONFERENCES/WORKSHOPS/SEMINAR ...\" - 3 columns\n// (blank spacer, label, citation). The single entry \"Mouli Rastogi ...\"\n// (originally row 2) is removed, so every later row's citation text\n// shifts up by one and the trailing row is deleted.\n// ---------------------------------------------------------------------\nconst confTable = tables.items[2];\nconfTable.rows.load(\"items\");\nawait context.sync();\n\nconst confRows = confTable.rows.items;\nfor (let r = 0; r < confRows.length; r++) {\n  confRows[r].cells.load(\"items\");\n}\nawait context.sync();\n\nfor (let r = 0; r < confRows.length; r++) {\n  const cells = confRows[r].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    cells[c].body.paragraphs.load(\"items\");\n  }\n}\nawait context.sync();\n\nconst confCitations = [];\nfor (let r = 0; r < confRows.length; r++) {\n  const cells = confRows[r].cells.items;\n  const lastCellParas = cells[cells.length - 1].body.paragraphs.items;\n  lastCellParas[0].load(\"text\");\n}\nawait context.sync();\nfor (let r = 0; r < confRows.length; r++) {\n  const cells = confRows[r].cells.items;\n  const lastCellParas = cells[cells.length - 1].body.paragraphs.items;\n  confCitations.push(lastCellParas[0].text);\n}\n\nconst removedIndex = 1; // 0-based index of the row being removed (\"Mouli Rastogi ...\", originally row 2)\n\nfor (let r = removedIndex; r < confRows.length - 1; r++) {\n  const cells = confRows[r].cells.items;\n  const lastCellParas = cells[cells.length - 1].body.paragraphs.items;\n  lastCellParas[0].insertText(confCitations[r + 1], Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconfTable.rows.load(\"items\");\nawait context.sync();\nconst confRowsAfter = confTable.rows.items;\nconfRowsAfter[confRowsAfter.length - 1].delete();\nawait context.sync();\n", "ps1": "# The document has two \"numbered\" tables whose first (small) column(s)\n# hold a static \"N.\" label and whose last column holds the publication /\n# conference citation text. This edit removes specific entries from each\n# table. Because the \"N.\" labels are plain static text (not an automatic\n# Word numbering field), removing an entry in the middle is implemented by\n# shifting the citation text of every later row up by one slot and then\n# physically deleting the now-duplicated row(s) at the end of the table -\n# this keeps the \"1.\", \"2.\", \"3.\" ... labels correct without having to\n# rewrite them.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Table 1: \"Papers Published by Faculty in Journals\" - 2 columns\n# (label, citation). The first two entries are removed, so the citation\n# text that used to belong to rows 3/4/5 becomes the new rows 1/2/3, and\n# the trailing two rows are deleted.\n# ---------------------------------------------------------------------\n$journalTable = $d.Tables.Item(1)\n\n$journalRowCount = $journalTable.Rows.Count\n$journalLastCol = $journalTable.Columns.Count\n$journalCitations = @()\nfor ($i = 1; $i -le $journalRowCount; $i++) {\n    $journalCitations += $journalTable.Rows.Item($i).Cells.Item($journalLastCol).Range.Text\n}\n\n$journalEntriesRemoved = 2\nfor ($i = 1; $i -le ($journalRowCount - $journalEntriesRemoved); $i++) {\n    $journalTable.Rows.Item($i).Cells.Item($journalLastCol).Range.Text = $journalCitations[$i - 1 + $journalEntriesRemoved]\n}\n\nfor ($i = 1; $i -le $journalEntriesRemoved; $i++) {\n    $journalTable.Rows.Item($journalTable.Rows.Count).Delete()\n}\n\n# ---------------------------------------------------------------------\n# Table 2: \"CONFERENCES/WORKSHOPS/SEMINAR ...\" - 3 columns\n# (blank spacer, label, citation). The single entry \"Mouli Rastogi ...\"\n# (originally row 2) is removed, so every later row's citation text\n# shifts up by one and the trailing row is deleted.\n# ---------------------------------------------------------------------\n$confTable = $d.Tables.Item(3)\n\n$confRowCount = $confTable.Rows.Count\n$confLastCol = $confTable.Columns.Count\n$confCitations = @()\nfor ($i = 1; $i -le $confRowCount; $i++) {\n    $confCitations += $confTable.Rows.Item($i).Cells.Item($confLastCol).Range.Text\n}\n\n$removedRow = 2\nfor ($i = $removedRow; $i -le ($confRowCount - 1); $i++) {\n    $confTable.Rows.Item($i).Cells.Item($confLastCol).Range.Text = $confCitations[$i]\n}\n\n$confTable.Rows.Item($confTable.Rows.Count).Delete()\n"}
